$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "-Function to check upvote status ... & display of most
#    upvoted questions/answers": remove the <w:proofErr> spell-check
#    wrapping around "upvote" / "upvoted" by merging each split run-triplet
#    into a single run, while leaving the middle " & " run untouched.
#
#    The engine renormalises (merges) every contiguous block of same-
#    format runs in a paragraph whenever it is edited, so a naive
#    Find/Replace would also swallow the " & " run. To avoid that we
#    temporarily fence the " & " run off with a scratch bookmark (bookmarks
#    act as a hard boundary for the left-hand side of a merge), do the
#    text surgery with an insert-before-delete pattern (which does not
#    trigger the "heal" merge that a delete-before-insert pattern does),
#    and finally drop the scratch bookmark again.
# ---------------------------------------------------------------------------

$pUpvote = $d.Paragraphs(9)
$pStart = $pUpvote.Range.Start

# Fence the " & " run ([51,54) relative to the paragraph) with a throwaway
# bookmark so it cannot be absorbed by the renormalisation pass.
$ampRange = $d.Range($pStart + 51, $pStart + 54)
$d.Bookmarks.Add("zzScratchAmp", $ampRange)

# Right-hand side: "display of most " + "upvoted" + " questions/answers"
# -> "display of most upvoted questions/answers"  (relative [54,95))
$insRight = $d.Range($pStart + 95, $pStart + 95)
$insRight.InsertBefore("display of most upvoted questions/answers")
$oldRight = $d.Range($pStart + 54, $pStart + 95)
$oldRight.Delete()

# Left-hand side: "-Function to check " + "upvote" + " status " + endash +
# " question/answers" -> single run (relative [0,51))
$insLeft = $d.Range($pStart + 51, $pStart + 51)
$insLeft.InsertBefore("-Function to check upvote status " + [char]0x2013 + " question/answers")
$oldLeft = $d.Range($pStart + 0, $pStart + 51)
$oldLeft.Delete()

# Remove the scratch bookmark; the " & " run it protected is left in place.
$d.Bookmarks("zzScratchAmp").Delete()

# ---------------------------------------------------------------------------
# 2) Paragraph "-Function to check for the ability to change user name":
#    move the run ahead of the _GoBack bookmark (bookmarkStart/bookmarkEnd
#    now follow the text instead of preceding it).
# ---------------------------------------------------------------------------

$pName = $d.Paragraphs(17)
$pNameEnd = $pName.Range.End

$d.Bookmarks("_GoBack").Delete()

# Insert a scratch character at the end of the paragraph's text so the
# bookmark can be anchored strictly after the run (anchoring a bookmark at
# the very last valid text offset of a paragraph misplaces it in this
# runtime), then delete the scratch character again.
$insScratch = $d.Range($pNameEnd - 1, $pNameEnd - 1)
$insScratch.InsertBefore("Z")

$pName2 = $d.Paragraphs(17)
$pName2End = $pName2.Range.End
$bmPos = $d.Range($pName2End - 2, $pName2End - 2)
$d.Bookmarks.Add("_GoBack", $bmPos)

$pName3 = $d.Paragraphs(17)
$pName3End = $pName3.Range.End
$scratchRange = $d.Range($pName3End - 2, $pName3End - 1)
$scratchRange.Delete()

# ---------------------------------------------------------------------------
# 3) Drop the trailing empty paragraphs / tab-stop paragraph that followed
#    the "change user name" paragraph, right before the section break.
# ---------------------------------------------------------------------------

$pFirstTail = $d.Paragraphs(18)
$pLastTail = $d.Paragraphs(21)
$tailRange = $d.Range($pFirstTail.Range.Start, $pLastTail.Range.End)
$tailRange.Delete()
